$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new competitor row (row 37) at the bottom of the table.
# Values are set in the same order as the original authoring so that the
# shared-strings table grows with matching order/indices.
$ws.Range("A37").Value = "Any.DO"
$ws.Range("B37").Value = "To do list"
$ws.Range("J37").Value = "downloaded 0.5M in 30days: http://techcrunch.com/2011/12/12/any-do-android-500000/"
$ws.Range("C37").Value = "http://www.any.do/"
$ws.Range("D37").Value = "GA"
$ws.Range("H37").Value = "android (iPhone, web coming)"

# Grow the worksheet's table so the new row is included in it.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:J37")) | Out-Null

# Match the saved selection/active cell on the sheet.
$ws.Range("C37").Select() | Out-Null
